# Updates the cryptos price/volume table with the latest scraped values.
# Column layout: A=rank(unchanged) B=Coin C=Link D=Price E=Volume(1h)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper that writes a value as plain text (avoids Excel auto-converting
# strings like "214.22" or "1.00" into numbers) and then resets the cell
# style back to the sheet's normal/unstyled look (matching column B of the
# same row, which never carries an explicit style) so no stray formatting
# is introduced by the text coercion.
function Set-TextValue($rowNum, $colNum, $val) {
    $cell = $ws.Cells.Item($rowNum, $colNum)
    $cell.Value = "'" + $val
    $cell.Style = $ws.Cells.Item($rowNum, 2).Style
}

Set-TextValue 2 4 '25.704.60'
Set-TextValue 2 5 '  -0.73%  '
Set-TextValue 3 4 '1.626.57'
Set-TextValue 3 5 '  -1.12%  '
Set-TextValue 4 4 '1.00'
Set-TextValue 4 5 '  -0.21%  '
Set-TextValue 5 4 '214.22'
Set-TextValue 5 5 '  -0.88%  '
Set-TextValue 6 4 '0.500'
Set-TextValue 6 5 '  -1.19%  '
Set-TextValue 7 5 '  -0.06%  '
Set-TextValue 8 5 '  -0.76%  '
Set-TextValue 9 4 '0.0635'
Set-TextValue 9 5 '  -1.31%  '
Set-TextValue 10 4 '19.55'
Set-TextValue 10 5 '  -4.21%  '
Set-TextValue 11 4 '0.0783'
Set-TextValue 11 5 '  +0.27%  '
Set-TextValue 12 4 '1.643.69'
Set-TextValue 12 5 '  -0.04%  '
Set-TextValue 13 4 '4.23'
Set-TextValue 13 5 '  -1.02%  '
Set-TextValue 14 4 '1.849.57'
Set-TextValue 14 5 '  -1.16%  '
Set-TextValue 15 4 '0.551'
Set-TextValue 15 5 '  -2.12%  '
Set-TextValue 16 4 '0.0₃0762'
Set-TextValue 16 5 '  -1.07%  '
Set-TextValue 17 4 '62.63'
Set-TextValue 17 5 '  -1.16%  '
Set-TextValue 18 4 '25.691.26'
Set-TextValue 18 5 '  -0.85%  '
Set-TextValue 19 5 '  -0.05%  '
Set-TextValue 20 4 '4.41'
Set-TextValue 20 5 '  +0.88%  '
Set-TextValue 21 4 '193.19'
Set-TextValue 21 5 '  +0.18%  '
Set-TextValue 22 4 '9.91'
Set-TextValue 22 5 '  -0.32%  '
Set-TextValue 23 4 '6.18'
Set-TextValue 23 5 '  +0.92%  '
Set-TextValue 24 5 '  +0.17%  '
Set-TextValue 25 4 '1.78'
Set-TextValue 25 5 '  -1.18%  '
Set-TextValue 26 4 '139.69'
Set-TextValue 26 5 '  -1.37%  '
Set-TextValue 27 4 '0.120'
Set-TextValue 27 5 '  -3.25%  '
Set-TextValue 28 4 '6.81'
Set-TextValue 28 5 '  +0.27%  '
Set-TextValue 29 4 '15.45'
Set-TextValue 29 5 '  -0.65%  '
Set-TextValue 30 5 '  -0.92%  '
Set-TextValue 31 4 '0.0485'
Set-TextValue 31 5 '  -1.91%  '
Set-TextValue 32 4 '3.31'
Set-TextValue 32 5 '  +0.30%  '
Set-TextValue 33 4 '3.23'
Set-TextValue 33 5 '  -0.23%  '
Set-TextValue 34 5 '  +0.37%  '
Set-TextValue 35 5 '  -0.11%  '
Set-TextValue 36 4 '0.893'
Set-TextValue 36 5 '  -1.45%  '
Set-TextValue 37 5 '  +0.21%  '
Set-TextValue 38 4 '0.543'
Set-TextValue 38 5 '  -2.43%  '
Set-TextValue 39 4 '1.106.59'
Set-TextValue 39 5 '  -2.32%  '
Set-TextValue 40 4 '0.0155'
Set-TextValue 40 5 '  -0.94%  '
Set-TextValue 41 5 '  +0.05%  '
Set-TextValue 42 5 '  +0.68%  '
Set-TextValue 43 4 '99.95'
Set-TextValue 43 5 '  +0.95%  '
Set-TextValue 44 4 '0.795'
Set-TextValue 44 5 '  -1.18%  '
Set-TextValue 45 4 '1.757.83'
Set-TextValue 46 5 '  -2.51%  '
Set-TextValue 47 4 '54.88'
Set-TextValue 47 5 '  -1.66%  '
Set-TextValue 48 5 '  -2.70%  '
Set-TextValue 49 4 '2.38'
Set-TextValue 49 5 '  +4.60%  '
Set-TextValue 50 2 'EnergySwap'
Set-TextValue 50 3 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 50 4 '7.65'
Set-TextValue 50 5 '  -1.18%  '
Set-TextValue 51 2 'Cronos'
Set-TextValue 51 3 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue 51 4 '0.0501'
Set-TextValue 51 5 '  -0.83%  '
